$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 4436
$ws.Range("J40").Value = 5760.6
$ws.Range("L40").Value = 5760.6
$ws.Range("N40").Value = -6110.6
$ws.Range("H53").Value = 1025.875
$ws.Range("I53").Value = 597.4
$ws.Range("J53").Value = 1740
$ws.Range("K53").Value = 597.4
$ws.Range("L53").Value = 1740
$ws.Range("M53").Value = 39.60000000000002
$ws.Range("N53").Value = -3014
$ws.Range("H96").Value = 263.06668
$ws.Range("I96").Value = 176.61539
$ws.Range("K96").Value = 529.84617
$ws.Range("M96").Value = 843.15383
$ws.Range("H97").Value = 911.2
$ws.Range("J97").Value = 911.2
$ws.Range("L97").Value = 2733.6
$ws.Range("N97").Value = -3725.6
$ws.Range("H100").Value = 4998.6665
$ws.Range("I100").Value = 3500.6667
$ws.Range("J100").Value = 5747.6665
$ws.Range("K100").Value = 3500.6667
$ws.Range("L100").Value = 5747.6665
$ws.Range("M100").Value = -2959.6667
$ws.Range("N100").Value = -6829.6665
$ws.Range("H104").Value = 1092
$ws.Range("J104").Value = 1089
$ws.Range("L104").Value = 3267
$ws.Range("N104").Value = -6761
$ws.Range("H132").Value = 14258.053
$ws.Range("I132").Value = 13656.4375
$ws.Range("K132").Value = 40969.3125
$ws.Range("M132").Value = -38439.3125
$ws.Range("H138").Value = 2117.2856
$ws.Range("I138").Value = 2095.2727
$ws.Range("J138").Value = 2198
$ws.Range("K138").Value = 6285.8181
$ws.Range("L138").Value = 6594
$ws.Range("M138").Value = -1145.8181
$ws.Range("N138").Value = -16874
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3476.5
$ws.Range("I45").Value = 2537.8572
$ws.Range("K45").Value = 2537.8572
$ws.Range("M45").Value = -2160.8572
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 716.5
$ws.Range("I29").Value = 716.5
$ws.Range("K29").Value = 716.5
$ws.Range("M29").Value = -427.5
$ws.Range("H37").Value = 2527
$ws.Range("I37").Value = 2602.2856
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 2602.2856
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = -2465.2856
$ws.Range("N37").Value = -2274
$ws.Range("H94").Value = 1309.6
$ws.Range("I94").Value = 1309.6
$ws.Range("K94").Value = 1309.6
$ws.Range("M94").Value = -858.5999999999999
$ws.Range("H99").Value = 3201.0557
$ws.Range("I99").Value = 2752.5386
$ws.Range("K99").Value = 2752.5386
$ws.Range("M99").Value = -1254.5386
$ws.Range("H134").Value = 2356.4707
$ws.Range("I134").Value = 2356.4707
$ws.Range("K134").Value = 7069.4121
$ws.Range("M134").Value = -4534.4121
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 39999.5
$ws.Range("J92").Value = 39999.5
$ws.Range("L92").Value = 39999.5
$ws.Range("N92").Value = -44991.5
$ws.Range("H105").Value = 1004.5
$ws.Range("I105").Value = 1035.8572
$ws.Range("J105").Value = 931.3333
$ws.Range("K105").Value = 1035.8572
$ws.Range("L105").Value = 931.3333
$ws.Range("M105").Value = 711.1428000000001
$ws.Range("N105").Value = -4425.3333
$ws.Range("H125").Value = 76542
$ws.Range("J125").Value = 76542
$ws.Range("L125").Value = 76542
$ws.Range("N125").Value = -81462
$ws.Range("H132").Value = 4674.0557
$ws.Range("J132").Value = 5766.4
$ws.Range("L132").Value = 17299.2
$ws.Range("N132").Value = -22359.2
$ws.Range("H134").Value = 2785.0715
$ws.Range("I134").Value = 2749.5
$ws.Range("K134").Value = 8248.5
$ws.Range("M134").Value = -5713.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10308850
$ws.Range("I11").Value = 11579368
$ws.Range("J11").Value = 6860300
$ws.Range("K11").Value = 11579368
$ws.Range("L11").Value = 6860300
$ws.Range("M11").Value = -11579229
$ws.Range("N11").Value = -6860578
$ws.Range("H80").Value = 2786.5
$ws.Range("I80").Value = 2505
$ws.Range("K80").Value = 2505
$ws.Range("M80").Value = -1507
$ws.Range("H83").Value = 2786.5
$ws.Range("I83").Value = 2505
$ws.Range("K83").Value = 12525
$ws.Range("M83").Value = -7533
$ws.Range("H92").Value = 5218.5
$ws.Range("J92").Value = 5821.143
$ws.Range("L92").Value = 5821.143
$ws.Range("N92").Value = -9565.143
$ws.Range("H97").Value = 579.64703
$ws.Range("I97").Value = 553.4375
$ws.Range("K97").Value = 553.4375
$ws.Range("M97").Value = -57.4375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1574.5
$ws.Range("J22").Value = 550
$ws.Range("L22").Value = 550
$ws.Range("N22").Value = -1140
$ws.Range("H27").Value = 1574.5
$ws.Range("J27").Value = 550
$ws.Range("L27").Value = 550
$ws.Range("N27").Value = -764
$ws.Range("H31").Value = 951.4286
$ws.Range("I31").Value = 566.25
$ws.Range("J31").Value = 1465
$ws.Range("K31").Value = 566.25
$ws.Range("L31").Value = 1465
$ws.Range("M31").Value = -318.25
$ws.Range("N31").Value = -1961
$ws.Range("H40").Value = 5742.6665
$ws.Range("I40").Value = 4778.4
$ws.Range("J40").Value = 6948
$ws.Range("K40").Value = 4778.4
$ws.Range("L40").Value = 6948
$ws.Range("M40").Value = -4642.4
$ws.Range("N40").Value = -7220
$ws.Range("H93").Value = 1185.6428
$ws.Range("I93").Value = 1202.4166
$ws.Range("K93").Value = 1202.4166
$ws.Range("M93").Value = 45.58339999999998
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H136").Value = 2333.1667
$ws.Range("I136").Value = 1299.8
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 3899.4
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -1349.4
$ws.Range("N136").Value = -27600
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 11851.167
$ws.Range("I3").Value = 3034.3333
$ws.Range("K3").Value = 3034.3333
$ws.Range("M3").Value = -2920.3333
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H96").Value = 2500.4
$ws.Range("I96").Value = 1875.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1875.5
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -502.5
$ws.Range("N96").Value = -7746
